$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range('A2').Value = 'Última actualización: 14:24:30'
$ws.Range('A3').Value = 'Total filas: 218'

$ws.Range('C68').Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range('C69').Value = '215A_EL PATO'
$ws.Range('C70').Value = '11_ETCHEVERRY'
$ws.Range('C76').Value = '16_SANTA ANA'
$ws.Range('C77').Value = '16_P MOR-SANTA ANA'
$ws.Range('C116').Value = '16_P MOR-167 Y 521'
$ws.Range('C117').Value = '11_ETCHEVERRY'
$ws.Range('C118').Value = '10_OLMOS'
$ws.Range('C119').Value = '23_HERNANDEZ'
$ws.Range('C142').Value = '16_P MOR-SANTA ANA'
$ws.Range('C143').Value = '23_HERNANDEZ'
$ws.Range('C153').Value = '11_ETCHEVERRY'
$ws.Range('C154').Value = '23_HERNANDEZ'
$ws.Range('C163').Value = '15_ABASTO'
$ws.Range('C164').Value = '14_ABASTO'
$ws.Range('C179').Value = '16_SANTA ANA'
$ws.Range('C180').Value = '17_ROMERO'
$ws.Range('C181').Value = '11_ETCHEVERRY'
$ws.Range('C182').Value = '215A_EL PATO'
$ws.Range('A184').Value = '12:59:25'
$ws.Range('C184').Value = '225_GOMEZ'
$ws.Range('D184').Value = 57
$ws.Range('A185').Value = '13:45:48'
$ws.Range('C185').Value = '16_P MOR-167 Y 521'
$ws.Range('D185').Value = 11
$ws.Range('C186').Value = '17_ROMERO'
$ws.Range('C187').Value = '23_HERNANDEZ'
$ws.Range('A195').Value = '14:24:30'
$ws.Range('B195').Value = '14:26'
$ws.Range('C195').Value = '16_SANTA ANA'
$ws.Range('D195').Value = 2
$ws.Range('A196').Value = '14:24:30'
$ws.Range('B196').Value = '14:28'
$ws.Range('C196').Value = '15_ABASTO'
$ws.Range('D196').Value = 4
$ws.Range('A197').Value = '14:24:30'
$ws.Range('B197').Value = '14:34'
$ws.Range('C197').Value = '23_HERNANDEZ'
$ws.Range('D197').Value = 10
$ws.Range('A198').Value = '12:59:25'
$ws.Range('B198').Value = '14:42'
$ws.Range('C198').Value = '14_ABASTO'
$ws.Range('D198').Value = 103
$ws.Range('A199').Value = '14:24:30'
$ws.Range('B199').Value = '14:44'
$ws.Range('C199').Value = '14_ABASTO'
$ws.Range('D199').Value = 20
$ws.Range('A200').Value = '14:24:30'
$ws.Range('B200').Value = '14:46'
$ws.Range('C200').Value = '16_SANTA ANA'
$ws.Range('D200').Value = 22
$ws.Range('A201').Value = '14:24:30'
$ws.Range('B201').Value = '14:56'
$ws.Range('C201').Value = '16_P MOR-SANTA ANA'
$ws.Range('D201').Value = 32
$ws.Range('A202').Value = '14:24:30'
$ws.Range('B202').Value = '14:58'
$ws.Range('C202').Value = '215B_EL PATO'
$ws.Range('D202').Value = 34
$ws.Range('A203').Value = '14:24:30'
$ws.Range('B203').Value = '15:00'
$ws.Range('C203').Value = '81_EL PELIGRO'
$ws.Range('D203').Value = 36
$ws.Range('A204').Value = '14:24:30'
$ws.Range('B204').Value = '15:05'
$ws.Range('C204').Value = '10_OLMOS'
$ws.Range('D204').Value = 41
$ws.Range('A205').Value = '14:24:30'
$ws.Range('B205').Value = '15:10'
$ws.Range('C205').Value = '17_ROMERO'
$ws.Range('D205').Value = 46
$ws.Range('A206').Value = '14:24:30'
$ws.Range('B206').Value = '15:13'
$ws.Range('C206').Value = '11_ETCHEVERRY'
$ws.Range('D206').Value = 49
$ws.Range('A207').Value = '14:24:30'
$ws.Range('B207').Value = '15:20'
$ws.Range('C207').Value = '15_ABASTO'
$ws.Range('D207').Value = 56
$ws.Range('A208').Value = '14:24:30'
$ws.Range('B208').Value = '15:21'
$ws.Range('C208').Value = '26_HERNANDEZ'
$ws.Range('D208').Value = 57
$ws.Range('A209').Value = '13:45:48'
$ws.Range('B209').Value = '15:22'
$ws.Range('C209').Value = '26_HERNANDEZ'
$ws.Range('D209').Value = 97
$ws.Range('E209').Value = 'LP1912'
$ws.Range('A210').Value = '14:24:30'
$ws.Range('B210').Value = '15:32'
$ws.Range('C210').Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range('D210').Value = 68
$ws.Range('E210').Value = 'LP1912'
$ws.Range('A211').Value = '13:45:48'
$ws.Range('B211').Value = '15:34'
$ws.Range('C211').Value = '23_HERNANDEZ'
$ws.Range('D211').Value = 109
$ws.Range('E211').Value = 'LP1912'
$ws.Range('A212').Value = '14:24:30'
$ws.Range('B212').Value = '15:37'
$ws.Range('C212').Value = '10_OLMOS'
$ws.Range('D212').Value = 73
$ws.Range('E212').Value = 'LP1912'
$ws.Range('A213').Value = '14:24:30'
$ws.Range('B213').Value = '15:38'
$ws.Range('C213').Value = '23_HERNANDEZ'
$ws.Range('D213').Value = 74
$ws.Range('E213').Value = 'LP1912'
$ws.Range('A214').Value = '14:24:30'
$ws.Range('B214').Value = '15:38'
$ws.Range('C214').Value = '215A_EL PATO'
$ws.Range('D214').Value = 74
$ws.Range('E214').Value = 'LP1912'
$ws.Range('A215').Value = '13:45:48'
$ws.Range('B215').Value = '15:42'
$ws.Range('C215').Value = '14_ABASTO'
$ws.Range('D215').Value = 117
$ws.Range('E215').Value = 'LP1912'
$ws.Range('A216').Value = '14:24:30'
$ws.Range('B216').Value = '15:45'
$ws.Range('C216').Value = '14_ABASTO'
$ws.Range('D216').Value = 81
$ws.Range('E216').Value = 'LP1912'
$ws.Range('A217').Value = '14:24:30'
$ws.Range('B217').Value = '15:46'
$ws.Range('C217').Value = '16_P MOR-167 Y 521'
$ws.Range('D217').Value = 82
$ws.Range('E217').Value = 'LP1912'
$ws.Range('A218').Value = '14:24:30'
$ws.Range('B218').Value = '15:53'
$ws.Range('C218').Value = '11_ETCHEVERRY'
$ws.Range('D218').Value = 89
$ws.Range('E218').Value = 'LP1912'
$ws.Range('A219').Value = '14:24:30'
$ws.Range('B219').Value = '15:56'
$ws.Range('C219').Value = '17_ROMERO'
$ws.Range('D219').Value = 92
$ws.Range('E219').Value = 'LP1912'
$ws.Range('A220').Value = '14:24:30'
$ws.Range('B220').Value = '15:56'
$ws.Range('C220').Value = '27_EL RETIRO'
$ws.Range('D220').Value = 92
$ws.Range('E220').Value = 'LP1912'
$ws.Range('A221').Value = '14:24:30'
$ws.Range('B221').Value = '16:15'
$ws.Range('C221').Value = '225_C ROCA-H SUR'
$ws.Range('D221').Value = 111
$ws.Range('E221').Value = 'LP1912'
$ws.Range('A222').Value = '14:24:30'
$ws.Range('B222').Value = '16:20'
$ws.Range('C222').Value = '215C_EL PATO'
$ws.Range('D222').Value = 116
$ws.Range('E222').Value = 'LP1912'
$ws.Range('A223').Value = '14:24:30'
$ws.Range('B223').Value = '16:21'
$ws.Range('C223').Value = '26_HERNANDEZ'
$ws.Range('D223').Value = 117
$ws.Range('E223').Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range('A2').Value = 'Última actualización: 14:24:30'
$ws.Range('A3').Value = 'Total filas: 24'

$ws.Range('A27').Value = '14:24:30'
$ws.Range('D27').Value = 34
$ws.Range('A28').Value = '14:24:30'
$ws.Range('D28').Value = 74
$ws.Range('A29').Value = '14:24:30'
$ws.Range('B29').Value = '16:20'
$ws.Range('C29').Value = '215C_EL PATO'
$ws.Range('D29').Value = 116
$ws.Range('E29').Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range('A2').Value = 'Última actualización: 14:24:30'
$ws.Range('A3').Value = 'Total filas: 31'

$ws.Range('A34').Value = '14:24:30'
$ws.Range('D34').Value = 29
$ws.Range('A35').Value = '14:24:30'
$ws.Range('D35').Value = 70
$ws.Range('A36').Value = '14:24:30'
$ws.Range('B36').Value = '16:14'
$ws.Range('C36').Value = '215C_LA PLATA'
$ws.Range('D36').Value = 110
$ws.Range('E36').Value = 'L6203'
